# Correção e atualização 13/11
# Remove the 8 "L2L3" measurement rows (old rows 4-11) from the Planilha1
# worksheet. These rows had Localização="L2L3", MUX=25, Channel 1-8, and no
# Nivel value. Deleting them shifts all subsequent rows up by 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Delete the entire rows 4 through 11 (inclusive), which removes the
# "L2L3" entries and shifts everything below up.
$ws.Range("A4:D11").EntireRow.Delete() | Out-Null

# Leave the selection near where the user last worked after the edit.
$ws.Activate()
$ws.Range("E14").Select() | Out-Null
